$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 23:35"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1642477
$ws.Range("C4").Value = 21575
$ws.Range("D4").Value = 395947
$ws.Range("E4").Value = 1148997
$ws.Range("G4").Value = 1179
$ws.Range("H4").Value = 97533

# --- Rows 110/111: Mali overtakes Guinea Ecuatorial, swapping their order ---
# Row 110 now corresponds to Mali (updated data), row 111 to Guinea Ecuatorial (old row-110 data)
$ws.Range("A110").Value = "Mali"
$ws.Range("B110").Value = 969
$ws.Range("C110").Value = 22
$ws.Range("D110").Value = 560
$ws.Range("E110").Value = 347
$ws.Range("G110").Value = 2
$ws.Range("H110").Value = 62

$ws.Range("A111").Value = "Guinea Ecuatorial"
$ws.Range("B111").Value = 960
$ws.Range("C111").Value = 57
$ws.Range("D111").Value = 165
$ws.Range("E111").Value = 784
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 11

# --- Row 154: Birmania ---
$ws.Range("D154").Value = 116
$ws.Range("E154").Value = 77
